# Relatório Atualizado e Enviado Para Avaliação
# Applies the commit's changes to the workbook via the Excel COM object model.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Workbook-level metadata
# ---------------------------------------------------------------------------

# _FilterDatabase defined name: normalise the sheet reference to the quoted
# form ('Relatorio'!$A$1:$O$23) that Excel wrote out after the edit.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Relatorio'!`$A`$1:`$O`$23"
    }
}

# ---------------------------------------------------------------------------
# 2) "Relatorio" worksheet - status/owner updates for rows 2, 3, 8, 9, 10, 11
# ---------------------------------------------------------------------------

$ws = $wb.Worksheets.Item("Relatorio")

# Rows 2 and 3 move from "Atribuído" to "Corrigido", and now record
# "Mateus Victor" as the corrector of the error (column O).
$ws.Range("I2").Value2 = "Corrigido"
$ws.Range("O2").Value2 = "Mateus Victor"

$ws.Range("I3").Value2 = "Corrigido"
$ws.Range("O3").Value2 = "Mateus Victor"

# Rows 8, 9 and 10 move from "Atribuído" to "Corrigido", now recording
# "Jeiel" as the corrector of the error (column O).
$ws.Range("I8").Value2 = "Corrigido"
$ws.Range("O8").Value2 = "Jeiel"

$ws.Range("I9").Value2 = "Corrigido"
$ws.Range("O9").Value2 = "Jeiel"

$ws.Range("I10").Value2 = "Corrigido"
$ws.Range("O10").Value2 = "Jeiel"

# Row 11 also moves from "Atribuído" to "Corrigido" (column O already held a
# value and stays untouched).
$ws.Range("I11").Value2 = "Corrigido"

# Rows 32-34: the "Titular do Erro" cell (column N) becomes horizontally
# centered, matching the style used elsewhere in the sheet.
$ws.Range("N32").HorizontalAlignment = -4108   # xlCenter
$ws.Range("N33").HorizontalAlignment = -4108   # xlCenter
$ws.Range("N34").HorizontalAlignment = -4108   # xlCenter

# Sheet view: scroll so column H is the left-most visible column, and move
# the active selection to O33.
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollColumn = 8
    $excel.ActiveWindow.ScrollRow = 1
} catch {
}
$ws.Range("O33").Select()

# ---------------------------------------------------------------------------
# 3) "grafico" worksheet - requalify the COUNTIF formulas to reference the
#    quoted sheet name ('Relatorio'!...) used after the edit. The cached
#    results (<v>) are refreshed automatically on recalculation.
# ---------------------------------------------------------------------------

$gws = $wb.Worksheets.Item("grafico")

$gws.Range("B2").Formula = '=COUNTIF(''Relatorio''!I2:I65,"Corrigido")'
$gws.Range("G2").Formula = '=COUNTIF(''Relatorio''!F2:F65,"Crítico")'

$gws.Range("B3").Formula = '=COUNTIF(''Relatorio''!I2:I65,"Atribuído")'
$gws.Range("G3").Formula = '=COUNTIF(''Relatorio''!F2:F65,"Moderado")'

$gws.Range("B4").Formula = '=COUNTIF(''Relatorio''!I2:I65,"Não Atribuído")'
$gws.Range("G4").Formula = '=COUNTIF(''Relatorio''!F2:F65,"Leve")'

$gws.Range("B5").Formula = '=COUNTIF(''Relatorio''!I5:I68,"Sem Ação")'

# Recalculate so every formula (and the chart series that reads from this
# sheet) reflects the updated "Relatorio" data.
$excel.CalculateFull()
